# Applies the bilibili-scraped 杭州-漫展信息 refresh described in the commit
# "Update gh-pages to output generated at 456a3b4": refreshed 想去人数 (want-to-go)
# counts across all sheets, dropped the finished 2024-08-04 SK-LiLi... wait event row,
# shifted the remaining August rows up by one, and appended the new 浮游猫动漫嘉年华 row.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions): refresh F-column (want-to-go counts) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Cells.Item(2, 6).Value = 115
$ws1.Cells.Item(3, 6).Value = 1044
$ws1.Cells.Item(4, 6).Value = 9248
$ws1.Cells.Item(5, 6).Value = 194
$ws1.Cells.Item(7, 6).Value = 1958
$ws1.Cells.Item(8, 6).Value = 6387
$ws1.Cells.Item(9, 6).Value = 619
$ws1.Cells.Item(11, 6).Value = 9744
$ws1.Cells.Item(12, 6).Value = 11061
$ws1.Cells.Item(13, 6).Value = 1233
$ws1.Cells.Item(14, 6).Value = 1140
$ws1.Cells.Item(15, 6).Value = 4910
$ws1.Cells.Item(16, 6).Value = 790
$ws1.Cells.Item(17, 6).Value = 447
$ws1.Cells.Item(19, 6).Value = 329
$ws1.Cells.Item(22, 6).Value = 237
$ws1.Cells.Item(23, 6).Value = 877
$ws1.Cells.Item(24, 6).Value = 1231
$ws1.Cells.Item(27, 6).Value = 2019
$ws1.Cells.Item(29, 6).Value = 614
$ws1.Cells.Item(30, 6).Value = 2649
$ws1.Cells.Item(32, 6).Value = 183
$ws1.Cells.Item(33, 6).Value = 1735

# --- Sheet "展览": rows 35-48 shift up one (old row 35 event concluded) and a new
# row 48 (杭州·浮游猫动漫嘉年华) is appended before the unaffected closing row 49 ---
$ws1.Cells.Item(35, 2).Value = '''2024-08-10'
$ws1.Cells.Item(35, 3).Value = '杭州·SK-LiLi综合同人展 × KK WORLD MINI快看漫画乐园'
$ws1.Cells.Item(35, 4).Value = '鸿泰路与明月桥路交汇处东南角方位(杭港地铁1号线/杭州地铁4号线彭埠站D口20米) 港龙悠乐城'
$ws1.Cells.Item(35, 5).Value = '2024.08.10 10:00-08.11 17:00'
$ws1.Cells.Item(35, 6).Value = 793
$ws1.Cells.Item(35, 7).Value = 68
$ws1.Cells.Item(35, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=88062'
$ws1.Cells.Item(35, 9).Value = '//i2.hdslb.com/bfs/openplatform/202406/xlLi5AXx1719324778808.jpeg'

$ws1.Cells.Item(36, 2).Value = '''2024-08-10'
$ws1.Cells.Item(36, 3).Value = '杭州·原神X崩坏X星铁旅行盛宴·首展'
$ws1.Cells.Item(36, 4).Value = '德胜东路2539号 梦马汽车小镇'
$ws1.Cells.Item(36, 5).Value = '2024.08.10 10:00-08.11 17:00'
$ws1.Cells.Item(36, 6).Value = 45
$ws1.Cells.Item(36, 7).Value = 65
$ws1.Cells.Item(36, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=88429'
$ws1.Cells.Item(36, 9).Value = '//i1.hdslb.com/bfs/openplatform/202407/1oKYAwlD1719816495200.png'

$ws1.Cells.Item(37, 2).Value = '''2024-08-10'
$ws1.Cells.Item(37, 3).Value = '杭州·原神X星铁X绝区零only'
$ws1.Cells.Item(37, 4).Value = '望江东路333号 杭州瑞莱克斯大酒店'
$ws1.Cells.Item(37, 5).Value = '2024.08.10 10:00-08.10 17:00'
$ws1.Cells.Item(37, 6).Value = 912
$ws1.Cells.Item(37, 7).Value = 60
$ws1.Cells.Item(37, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=82754'
$ws1.Cells.Item(37, 9).Value = '//i1.hdslb.com/bfs/openplatform/202403/qA0LNJuF1710234461030.jpeg'

$ws1.Cells.Item(38, 2).Value = '''2024-08-10'
$ws1.Cells.Item(38, 3).Value = '杭州·造梦探险家城堡二次元同好会'
$ws1.Cells.Item(38, 4).Value = '大岭山路156号 爱丽芬城堡'
$ws1.Cells.Item(38, 5).Value = '2024.08.10 10:00-08.10 22:00'
$ws1.Cells.Item(38, 6).Value = 589
$ws1.Cells.Item(38, 7).Value = 38
$ws1.Cells.Item(38, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=86432'
$ws1.Cells.Item(38, 9).Value = '//i2.hdslb.com/bfs/openplatform/202405/xWUy30Ns1716783723057.jpeg'

$ws1.Cells.Item(39, 2).Value = '''2024-08-10'
$ws1.Cells.Item(39, 3).Value = '杭州·首届LoveLive Only'
$ws1.Cells.Item(39, 4).Value = '同协路288号 1928创意园'
$ws1.Cells.Item(39, 5).Value = '2024.08.10 10:00-08.10 22:00'
$ws1.Cells.Item(39, 6).Value = 16
$ws1.Cells.Item(39, 7).Value = 28
$ws1.Cells.Item(39, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=88458'
$ws1.Cells.Item(39, 9).Value = '//i1.hdslb.com/bfs/openplatform/202406/cT9gANGo1719481033302.png'

$ws1.Cells.Item(40, 2).Value = '''2024-08-17'
$ws1.Cells.Item(40, 3).Value = '【会员购严选】杭州·首届次元之门动漫游戏博览会'
$ws1.Cells.Item(40, 4).Value = '阳城路雅澳杭州电商产业园西侧约200米 杭州大会展中心'
$ws1.Cells.Item(40, 5).Value = '2024.08.17 10:00-08.18 17:30'
$ws1.Cells.Item(40, 6).Value = 3304
$ws1.Cells.Item(40, 7).Value = 75
$ws1.Cells.Item(40, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=87065'
$ws1.Cells.Item(40, 9).Value = '//i1.hdslb.com/bfs/openplatform/202406/wrxORgrP1717593610187.jpeg'

$ws1.Cells.Item(41, 2).Value = '''2024-08-17'
$ws1.Cells.Item(41, 3).Value = '浙江·马娘ONLY03-晴风游憩'
$ws1.Cells.Item(41, 4).Value = '康候圣街99号 顺丰创新中心'
$ws1.Cells.Item(41, 5).Value = '2024.08.17 10:00-08.17 17:00'
$ws1.Cells.Item(41, 6).Value = 234
$ws1.Cells.Item(41, 7).Value = 68
$ws1.Cells.Item(41, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=86529'
$ws1.Cells.Item(41, 9).Value = '//i1.hdslb.com/bfs/openplatform/202405/21d6moub1716799089058.jpeg'

$ws1.Cells.Item(42, 2).Value = '''2024-08-18'
$ws1.Cells.Item(42, 3).Value = '杭州·少女番only3.0'
$ws1.Cells.Item(42, 4).Value = '凤起东路211号 名人名家宴会艺术中心(顺福店)'
$ws1.Cells.Item(42, 5).Value = '2024.08.18 10:00-08.18 17:00'
$ws1.Cells.Item(42, 6).Value = 81
$ws1.Cells.Item(42, 7).Value = 89
$ws1.Cells.Item(42, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=87676'
$ws1.Cells.Item(42, 9).Value = '//i1.hdslb.com/bfs/openplatform/202406/VHJXciCn1718433428129.jpeg'

$ws1.Cells.Item(43, 2).Value = '''2024-08-18'
$ws1.Cells.Item(43, 3).Value = '浙江·蔚蓝档案ONLY02-夏末狂欢！'
$ws1.Cells.Item(43, 4).Value = '康候圣街99号 顺丰创新中心'
$ws1.Cells.Item(43, 5).Value = '2024.08.18 10:00-08.18 17:00'
$ws1.Cells.Item(43, 6).Value = 505
$ws1.Cells.Item(43, 7).Value = 68
$ws1.Cells.Item(43, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=86594'
$ws1.Cells.Item(43, 9).Value = '//i1.hdslb.com/bfs/openplatform/202405/TVvJFURG1716799911888.jpeg'

$ws1.Cells.Item(44, 2).Value = '''2024-08-24'
$ws1.Cells.Item(44, 3).Value = '杭州·D3动漫游戏嘉年华'
$ws1.Cells.Item(44, 4).Value = '德胜东路2539号 梦马汽车小镇'
$ws1.Cells.Item(44, 5).Value = '2024.08.24 10:00-08.24 17:00'
$ws1.Cells.Item(44, 6).Value = 575
$ws1.Cells.Item(44, 7).Value = 50
$ws1.Cells.Item(44, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=84912'
$ws1.Cells.Item(44, 9).Value = '//i0.hdslb.com/bfs/openplatform/202406/LC5aplda1718697399055.jpeg'

$ws1.Cells.Item(45, 2).Value = '''2024-08-24'
$ws1.Cells.Item(45, 3).Value = '杭州·奶司的小人国娃展Nice Mini World'
$ws1.Cells.Item(45, 4).Value = '凤起东路211号 名人名家宴会艺术中心(顺福店)'
$ws1.Cells.Item(45, 5).Value = '2024.08.24 10:30-08.24 17:00'
$ws1.Cells.Item(45, 6).Value = 27
$ws1.Cells.Item(45, 7).Value = 60
$ws1.Cells.Item(45, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=86954'
$ws1.Cells.Item(45, 9).Value = '//i0.hdslb.com/bfs/openplatform/202406/aeSvo0X71717659096631.png'

$ws1.Cells.Item(46, 2).Value = '''2024-08-24'
$ws1.Cells.Item(46, 3).Value = '杭州·萌忧·原崩铁only'
$ws1.Cells.Item(46, 4).Value = '康候圣街99号 顺丰创新中心'
$ws1.Cells.Item(46, 5).Value = '2024.08.24 10:30-08.24 17:00'
$ws1.Cells.Item(46, 6).Value = 894
$ws1.Cells.Item(46, 7).Value = 50
$ws1.Cells.Item(46, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=87293'
$ws1.Cells.Item(46, 9).Value = '//i2.hdslb.com/bfs/openplatform/202406/rQFz5smR1717475284585.jpeg'

$ws1.Cells.Item(47, 2).Value = '''2024-08-24'
$ws1.Cells.Item(47, 3).Value = '杭州·金魂ONLY'
$ws1.Cells.Item(47, 4).Value = '金一路79号 XPACE湾区数字公园'
$ws1.Cells.Item(47, 5).Value = '2024.08.24 10:00-08.24 18:00'
$ws1.Cells.Item(47, 6).Value = 235
$ws1.Cells.Item(47, 7).Value = 78
$ws1.Cells.Item(47, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=87230'
$ws1.Cells.Item(47, 9).Value = '//i2.hdslb.com/bfs/openplatform/202406/SfNdoHfv1718176444920.png'

$ws1.Cells.Item(48, 2).Value = '''2024-09-15'
$ws1.Cells.Item(48, 3).Value = '杭州·浮游猫动漫嘉年华'
$ws1.Cells.Item(48, 4).Value = '鸿泰路与明月桥路交汇处东南角方位(杭港地铁1号线/杭州地铁4号线彭埠站D口20米) 港龙悠乐城'
$ws1.Cells.Item(48, 5).Value = '2024.09.15 09:00-09.16 18:00'
$ws1.Cells.Item(48, 6).Value = 3
$ws1.Cells.Item(48, 7).Value = 68
$ws1.Cells.Item(48, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=88498'
$ws1.Cells.Item(48, 9).Value = '//i2.hdslb.com/bfs/openplatform/202406/qsuFy4iv1719569431608.jpeg'

# Row 49 (杭州·理想乡动漫展-同人创作者大会) keeps its content; only the want-to-go count changes
$ws1.Cells.Item(49, 6).Value = 4200

# --- Sheet "本地生活" (Local life): refresh F2 want-to-go count ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Cells.Item(2, 6).Value = 5888

# --- Sheet "全部类型" (All types, merged view): refresh matching F-column counts ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Cells.Item(2, 6).Value = 1044
$ws4.Cells.Item(3, 6).Value = 9248
$ws4.Cells.Item(7, 6).Value = 6387
$ws4.Cells.Item(8, 6).Value = 619
$ws4.Cells.Item(9, 6).Value = 9744
$ws4.Cells.Item(10, 6).Value = 11061
$ws4.Cells.Item(12, 6).Value = 1233
$ws4.Cells.Item(13, 6).Value = 1140
$ws4.Cells.Item(14, 6).Value = 4910
$ws4.Cells.Item(15, 6).Value = 790
$ws4.Cells.Item(16, 6).Value = 447
$ws4.Cells.Item(21, 6).Value = 237
$ws4.Cells.Item(22, 6).Value = 877
$ws4.Cells.Item(23, 6).Value = 1231
$ws4.Cells.Item(26, 6).Value = 2019
$ws4.Cells.Item(28, 6).Value = 614
$ws4.Cells.Item(29, 6).Value = 2649
$ws4.Cells.Item(30, 6).Value = 183
$ws4.Cells.Item(31, 6).Value = 1735
$ws4.Cells.Item(33, 6).Value = 793
$ws4.Cells.Item(36, 6).Value = 45
$ws4.Cells.Item(37, 6).Value = 912
$ws4.Cells.Item(45, 6).Value = 575
